$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. First paragraph: append two trailing spaces to the existing run,
#    then append three new red (FF0000) runs reconstructing
#    "(This is a change – Version for main branch)" split exactly as
#    three separate runs, matching an incremental typing/edit session.
# -----------------------------------------------------------------

$p1 = $d.Paragraphs(1)
$origLen = $p1.Range.Text.Length - 1   # exclude trailing paragraph mark
$p1.Range.InsertAfter("  ")

$chunk1 = [string][char]0x0028 + "This is a change " + [string][char]0x2013 + " Ve"
$rA = $d.Paragraphs(1).Range
$rA.InsertAfter($chunk1)
$startA = $d.Paragraphs(1).Range.Start + $origLen + 2
$endA = $startA + $chunk1.Length
$d.Range($startA, $endA).Font.Color = 255

$chunk2 = "rsion for main branch"
$rB = $d.Paragraphs(1).Range
$beforeB = $rB.Text.Length
$rB.InsertAfter($chunk2)
$startB = $d.Paragraphs(1).Range.Start + ($beforeB - 1)
$endB = $startB + $chunk2.Length
$d.Range($startB, $endB).Font.Color = 255

$chunk3 = ")"
$rC = $d.Paragraphs(1).Range
$beforeC = $rC.Text.Length
$rC.InsertAfter($chunk3)
$startC = $d.Paragraphs(1).Range.Start + ($beforeC - 1)
$endC = $startC + $chunk3.Length
$d.Range($startC, $endC).Font.Color = 255

# -----------------------------------------------------------------
# 2. Remove the trailing "ank God almighty, we are free at last."
#    paragraph entirely (it followed "...Shall be lifted—nevermore!").
# -----------------------------------------------------------------

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$lastPara.Range.Delete()
